$wb = $excel.ActiveWorkbook

# Sheet ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6699.909
$ws.Range("I51").Value = 8670
$ws.Range("J51").Value = 5058.1665
$ws.Range("K51").Value = 8670
$ws.Range("L51").Value = 5058.1665
$ws.Range("M51").Value = -8186
$ws.Range("N51").Value = -6026.1665

# Sheet ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2629.7778
$ws.Range("I62").Value = 2629.7778
$ws.Range("K62").Value = 2629.7778
$ws.Range("M62").Value = -2005.7778

# Sheet ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2629.7778
$ws.Range("I65").Value = 2629.7778
$ws.Range("K65").Value = 13148.889
$ws.Range("M65").Value = -10028.889

# Sheet ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5761.25
$ws.Range("I69").Value = 5676.6665
$ws.Range("K69").Value = 17029.9995
$ws.Range("M69").Value = -16155.9995

# Sheet ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 5761.25
$ws.Range("I72").Value = 5676.6665
$ws.Range("K72").Value = 51089.9985
$ws.Range("M72").Value = -46721.9985

# Sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1304.3243
$ws.Range("I132").Value = 1093.1428
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 3279.4284
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -749.4284000000002
$ws.Range("N132").Value = -20060

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4917.826
$ws.Range("I137").Value = 5100.476
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 15301.428
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -12751.428
$ws.Range("N137").Value = -14100

# Sheet ALC row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 149980
$ws.Range("J139").Value = 149980
$ws.Range("L139").Value = 149980
$ws.Range("N139").Value = -160260

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10529.96
$ws.Range("I32").Value = 8011.2
$ws.Range("K32").Value = 8011.2
$ws.Range("M32").Value = -7724.2

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3108.2424
$ws.Range("I74").Value = 2344.0344
$ws.Range("J74").Value = 8648.75
$ws.Range("K74").Value = 2344.0344
$ws.Range("L74").Value = 8648.75
$ws.Range("M74").Value = -1470.0344
$ws.Range("N74").Value = -10396.75

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3108.2424
$ws.Range("I77").Value = 2344.0344
$ws.Range("J77").Value = 8648.75
$ws.Range("K77").Value = 11720.172
$ws.Range("L77").Value = 43243.75
$ws.Range("M77").Value = -7352.172
$ws.Range("N77").Value = -51979.75

# Sheet ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 840.9815
$ws.Range("I97").Value = 818.55554
$ws.Range("J97").Value = 953.1111
$ws.Range("K97").Value = 818.55554
$ws.Range("L97").Value = 953.1111
$ws.Range("M97").Value = -322.55554
$ws.Range("N97").Value = -1945.1111

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2536.4285
$ws.Range("I132").Value = 2536.4285
$ws.Range("K132").Value = 7609.2855
$ws.Range("M132").Value = -5079.2855

# Sheet BSM row 53
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = ""

# Sheet BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 43348
$ws.Range("J81").Value = 43348
$ws.Range("L81").Value = 43348
$ws.Range("N81").Value = -45470

# Sheet BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 43348
$ws.Range("J84").Value = 43348
$ws.Range("L84").Value = 130044
$ws.Range("N84").Value = -140652

# Sheet BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1310.5294
$ws.Range("J105").Value = 1505.5
$ws.Range("L105").Value = 1505.5
$ws.Range("N105").Value = -4999.5

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5486.1377
$ws.Range("I31").Value = 4134.0835
$ws.Range("K31").Value = 4134.0835
$ws.Range("M31").Value = -3839.0835

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5486.1377
$ws.Range("I34").Value = 4134.0835
$ws.Range("K34").Value = 4134.0835
$ws.Range("M34").Value = -3932.0835

# Sheet CRP row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 30980
$ws.Range("J68").Value = 27040
$ws.Range("L68").Value = 27040
$ws.Range("N68").Value = -28538

# Sheet CRP row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 30980
$ws.Range("J71").Value = 27040
$ws.Range("L71").Value = 81120
$ws.Range("N71").Value = -88608

# Sheet CRP row 97
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 27028.166
$ws.Range("J97").Value = 27028.166
$ws.Range("L97").Value = 27028.166
$ws.Range("N97").Value = -29010.166

# Sheet CRP row 102
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H102").Value = 195000
$ws.Range("I102").Value = 30000
$ws.Range("J102").Value = 360000
$ws.Range("K102").Value = 30000
$ws.Range("L102").Value = 360000
$ws.Range("M102").Value = -27566
$ws.Range("N102").Value = -364868

# Sheet CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 71429976
$ws.Range("J105").Value = 1849.5
$ws.Range("L105").Value = 1849.5
$ws.Range("N105").Value = -5343.5

# Sheet CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 112.4
$ws.Range("J23").Value = 138.125
$ws.Range("L23").Value = 414.375
$ws.Range("N23").Value = -884.375

# Sheet CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1304.3334
$ws.Range("I39").Value = 434.5
$ws.Range("J39").Value = 1552.8572
$ws.Range("K39").Value = 1303.5
$ws.Range("L39").Value = 4658.571599999999
$ws.Range("M39").Value = -1009.5
$ws.Range("N39").Value = -5246.571599999999

# Sheet CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 59782.883
$ws.Range("I121").Value = 360.44446
$ws.Range("K121").Value = 1081.33338
$ws.Range("M121").Value = 228.66662

# Sheet CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 529863.2
$ws.Range("I129").Value = 85499.836
$ws.Range("J129").Value = 1122347.6
$ws.Range("K129").Value = 256499.508
$ws.Range("L129").Value = 3367042.8
$ws.Range("M129").Value = -251499.508
$ws.Range("N129").Value = -3377042.8

# Sheet GSM row 51
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 55250
$ws.Range("J51").Value = 55250
$ws.Range("L51").Value = 55250
$ws.Range("N51").Value = -56268

# Sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1636.3334
$ws.Range("I102").Value = 1241.9744
$ws.Range("K102").Value = 1241.9744
$ws.Range("M102").Value = 380.0255999999999

# Sheet GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8559.875
$ws.Range("I126").Value = 5500
$ws.Range("K126").Value = 16500
$ws.Range("M126").Value = -14030

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7971.75
$ws.Range("I132").Value = 7404.4
$ws.Range("K132").Value = 22213.2
$ws.Range("M132").Value = -19683.2

# Sheet LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 125000170
$ws.Range("J55").Value = 223.33333
$ws.Range("L55").Value = 223.33333
$ws.Range("N55").Value = -569.3333299999999

# Sheet LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1940.3793
$ws.Range("I61").Value = 2207.3635
$ws.Range("K61").Value = 2207.3635
$ws.Range("M61").Value = -2005.3635

# Sheet LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5488.8
$ws.Range("I68").Value = 3814.6667
$ws.Range("K68").Value = 3814.6667
$ws.Range("M68").Value = -3065.6667

# Sheet LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5488.8
$ws.Range("I71").Value = 3814.6667
$ws.Range("K71").Value = 19073.3335
$ws.Range("M71").Value = -15329.3335

# Sheet LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1940.3793
$ws.Range("I113").Value = 2207.3635
$ws.Range("K113").Value = 2207.3635
$ws.Range("M113").Value = -37.36349999999993

# Sheet LTW row 138
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

# Sheet WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 76292.336
$ws.Range("J46").Value = 84495
$ws.Range("L46").Value = 84495
$ws.Range("N46").Value = -84957

# Sheet WVR row 70
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29105
$ws.Range("J70").Value = 29105
$ws.Range("L70").Value = 29105
$ws.Range("N70").Value = -29735

# Sheet WVR row 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 29105
$ws.Range("J73").Value = 29105
$ws.Range("L73").Value = 29105
$ws.Range("N73").Value = -31289

# Sheet WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2700
$ws.Range("I81").Value = 2700
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5400
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4339
$ws.Range("N81").Value = ""

# Sheet WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2700
$ws.Range("I84").Value = 2700
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 27000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -21696
$ws.Range("N84").Value = ""

# Sheet WVR row 95
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 57666
$ws.Range("I95").Value = 60000
$ws.Range("J95").Value = 56499
$ws.Range("K95").Value = 60000
$ws.Range("L95").Value = 56499
$ws.Range("M95").Value = -57254
$ws.Range("N95").Value = -61991

# Sheet WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 848.4737
$ws.Range("I113").Value = 639.1667
$ws.Range("J113").Value = 1207.2858
$ws.Range("K113").Value = 1917.5001
$ws.Range("L113").Value = 3621.8574
$ws.Range("M113").Value = 252.4999
$ws.Range("N113").Value = -7961.857400000001

# Sheet WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2387.5642
$ws.Range("I122").Value = 1591.7646
$ws.Range("K122").Value = 4775.293799999999
$ws.Range("M122").Value = -2325.293799999999

# Sheet WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 76292.336
$ws.Range("J134").Value = 84495
$ws.Range("L134").Value = 253485
$ws.Range("N134").Value = -258555

# Sheet WVR row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

# Sheet WVR row 140
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 93642.78
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 93642.78
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 93642.78
$ws.Range("M140").Value = ""
$ws.Range("N140").Value = -104002.78

# Sheet WVR row 141
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 83803.09
$ws.Range("J141").Value = 83803.09
$ws.Range("L141").Value = 83803.09
$ws.Range("N141").Value = -94163.09
